# Weekly update: insert a new Mercado Mayorista Lo Valledor de Santiago -
# Betarraga price block (Primera / Segunda / Tercera) for the newest date,
# pushing the existing history down by one block (3 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows right before the current row 427, shifting the
# rest of the table (and the old row 427-493 data) down to 430-496.
$ws.Range("A427:A429").EntireRow.Insert()

$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$catId = 100114014
$categoria = "Betarraga"
$variedad = "Sin especificar"
$unidad = "$/unidad"
$origen = "Región Metropolitana"
$kgUnidades = 1
$clasificacion = "Hortaliza"
$fecha = 44474

$newRows = @(
    @{ Row = 427; Calidad = "Primera"; Volumen = 49000; PMin = 110; PMax = 120; PProm = 115 },
    @{ Row = 428; Calidad = "Segunda"; Volumen = 44000; PMin = 90;  PMax = 95;  PProm = 92  },
    @{ Row = 429; Calidad = "Tercera"; Volumen = 15000; PMin = 70;  PMax = 70;  PProm = 70  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 6
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $catId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $kgUnidades
    $ws.Cells.Item($row, 18).Value = $clasificacion
}
